# Generate Report for Handoff
# Adds a new file "cc5f8fb9-7cfc-40ce-a5aa-1e91ff7cb0e4.md" (status: Ready for handoff)
# as a new row (row 9) on the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$guidFile  = "cc5f8fb9-7cfc-40ce-a5aa-1e91ff7cb0e4.md"
$e2ePath   = "e2e\cc5f8fb9-7cfc-40ce-a5aa-1e91ff7cb0e4.md"
$commitSha = "6353a72666a621ae894e13e975d5dfa86d8c9dcd"
$ghUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$guidFile"

$zhXlf = "cc5f8fb9-7cfc-40ce-a5aa-1e91ff7cb0e4.090878fb2b9a44d52ffe1853bfdde597708a66f6.zh-cn.xlf"
$deXlf = "cc5f8fb9-7cfc-40ce-a5aa-1e91ff7cb0e4.090878fb2b9a44d52ffe1853bfdde597708a66f6.de-de.xlf"

$zhHoDate = "2016-08-17 10:41:16"
$deHoDate = "2016-08-17 10:41:21"
$overviewDate = "2016-08-17 10:41:21"

$hyperlinkColor = 15570276  # BGR for RGB FF6495ED (cornflower blue), matches existing HyperLink style

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) - new row 9
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A9").Value = $guidFile
$wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), $ghUrl, "", "", $e2ePath)
$wsOverview.Range("B9").Font.Color = $hyperlinkColor
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = ""
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("G9").Value = $overviewDate

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) - new row 9
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A9"), $ghUrl, "", "", $guidFile)
$wsZh.Range("A9").Font.Color = $hyperlinkColor
$wsZh.Range("B9").Value = ".md"
$wsZh.Range("C9").Value = "Ready for handoff"
$wsZh.Range("D9").Value = "e2e"
$wsZh.Range("E9").Value = "ht"
$wsZh.Range("F9").Value = "'False"
$wsZh.Range("G9").Value = $zhXlf
$wsZh.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H9").Value = $zhHoDate
$wsZh.Range("I9").Value = ""
$wsZh.Range("J9").Value = ""
$wsZh.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K9").Value = "0001-01-01 00:00:00"
$wsZh.Range("L9").Value = ""
$wsZh.Range("M9").Value = "'True"
$wsZh.Range("N9").Value = ""
$wsZh.Range("O9").Value = "'False"
$wsZh.Range("P9").Value = ""

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) - new row 9
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A9"), $ghUrl, "", "", $guidFile)
$wsDe.Range("A9").Font.Color = $hyperlinkColor
$wsDe.Range("B9").Value = ".md"
$wsDe.Range("C9").Value = "Ready for handoff"
$wsDe.Range("D9").Value = "e2e"
$wsDe.Range("E9").Value = "ht"
$wsDe.Range("F9").Value = "'False"
$wsDe.Range("G9").Value = $deXlf
$wsDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H9").Value = $deHoDate
$wsDe.Range("I9").Value = ""
$wsDe.Range("J9").Value = ""
$wsDe.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDe.Range("L9").Value = ""
$wsDe.Range("M9").Value = "'True"
$wsDe.Range("N9").Value = ""
$wsDe.Range("O9").Value = "'False"
$wsDe.Range("P9").Value = ""
